$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.129.11'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '3.175.04'
$ws.Range("E3").Value = '  -4.63%  '
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '592.10'
$c.ClearFormats()
$ws.Range("E5").Value = '  -2.14%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '135.22'
$c.ClearFormats()
$ws.Range("E6").Value = '  -5.41%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.175.11'
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("E11").Value = '  -5.08%  '
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("E13").Value = '  -4.19%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '35.05'
$c.ClearFormats()
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '3.696.56'
$ws.Range("E15").Value = '  -4.66%  '
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").Value = '3.172.50'
$ws.Range("E17").Value = '  -4.53%  '
$ws.Range("D18").Value = '63.079.63'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("E19").Value = '  -4.06%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '462.50'
$c.ClearFormats()
$ws.Range("E20").Value = '  -3.93%  '
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("E22").Value = '  -5.14%  '
$ws.Range("E23").Value = '  -4.14%  '
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E29").Value = '  -6.65%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.81'
$c.ClearFormats()
$ws.Range("E30").Value = '  -5.48%  '
$ws.Range("E31").Value = '  -6.00%  '
$ws.Range("E32").Value = '  -6.02%  '
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("E34").Value = '  -6.52%  '
$ws.Range("E35").Value = '  -5.80%  '
$ws.Range("E36").Value = '  -4.05%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '51.47'
$c.ClearFormats()
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  -5.34%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0389'
$c.ClearFormats()
$ws.Range("E39").Value = '  -3.10%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '404.19'
$c.ClearFormats()
$ws.Range("E40").Value = '  -6.98%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '8.15'
$c.ClearFormats()
$ws.Range("E41").Value = '  -2.56%  '
$ws.Range("E42").Value = '  -4.63%  '
$ws.Range("D43").Value = '2.817.04'
$ws.Range("E43").Value = '  -9.09%  '
$ws.Range("E44").Value = '  -6.30%  '
$ws.Range("E45").Value = '  -5.89%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  -6.31%  '
$ws.Range("E48").Value = '  -4.13%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '123.83'
$c.ClearFormats()
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("E50").Value = '  -1.69%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '34.15'
$c.ClearFormats()
$ws.Range("E51").Value = '  -7.68%  '
